$wb = $excel.ActiveWorkbook

# --- Resultats sheet: enter game 7 (Q column) score for ST-LOUIS vs DALLAS (row 24/25) ---
$res = $wb.Worksheets.Item("Résultats")
$res.Range("Q24").Value = 4
$res.Range("Q25").Value = 3

# --- Pool sheet: swap the "Raymond Tiefengraber" (row 82) and "Michel Boulianne" (row 78)
#     entries: the name and each round's per-game score moved to the other row. ---
$pool = $wb.Worksheets.Item("Pool")

$cols = @("B","K","M","O","R","T","U","V","X","Y","Z","AA")
foreach ($col in $cols) {
    $cell78 = $pool.Range($col + "78")
    $cell82 = $pool.Range($col + "82")
    $v78 = $cell78.Value()
    $v82 = $cell82.Value()
    $cell78.Value = $v82
    $cell82.Value = $v78
}
